# -----------------------------------------------------------------------------
# PollsData.xlsx update: "update with cluster17 polls"
#
# 1. Poll-cluster 34 originally spanned rows 115-123; the last three of those
#    rows (the id column) are renumbered to poll-cluster 35.
# 2. Eight brand-new observation rows are appended (rows 124-131), covering
#    poll-clusters 36 (elabe), 37/38/39 (cluster17) - the new pollster named
#    in the commit message.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label rows 121-123 from poll-cluster 34 to poll-cluster 35 ----------
$ws.Cells.Item(121, 1).Value = 35
$ws.Cells.Item(122, 1).Value = 35
$ws.Cells.Item(123, 1).Value = 35

# --- row 124 ---
$ws.Cells.Item(124, 1).Value = 36
$ws.Cells.Item(124, 2).Value = 2021
$ws.Cells.Item(124, 3).Value = 14
$ws.Cells.Item(124, 4).Value = 11
$ws.Cells.Item(124, 5).Value = 24
$ws.Cells.Item(124, 6).Value = "elabe"
$ws.Cells.Item(124, 7).Value = "online"
$ws.Cells.Item(124, 8).Value = "partially"
$ws.Cells.Item(124, 9).Value = 876
$ws.Cells.Item(124, 10).Value = 1
$ws.Cells.Item(124, 11).Value = 1
$ws.Cells.Item(124, 12).Value = 9
$ws.Cells.Item(124, 13).Value = 2
$ws.Cells.Item(124, 14).Value = 1
$ws.Cells.Item(124, 15).Value = 8
$ws.Cells.Item(124, 16).Value = 4
$ws.Cells.Item(124, 17).Value = 25
$ws.Cells.Item(124, 20).Value = 13
$ws.Cells.Item(124, 21).Value = 1
$ws.Cells.Item(124, 22).Value = 2
$ws.Cells.Item(124, 23).Value = 20
$ws.Cells.Item(124, 24).Value = 12
$ws.Cells.Item(124, 25).Value = "T_1"
$ws.Cells.Item(124, 26).Value = "T_1"
$ws.Cells.Item(124, 27).Value = 1

# --- row 125 ---
$ws.Cells.Item(125, 1).Value = 36
$ws.Cells.Item(125, 2).Value = 2021
$ws.Cells.Item(125, 3).Value = 14
$ws.Cells.Item(125, 4).Value = 11
$ws.Cells.Item(125, 5).Value = 24
$ws.Cells.Item(125, 6).Value = "elabe"
$ws.Cells.Item(125, 7).Value = "online"
$ws.Cells.Item(125, 8).Value = "partially"
$ws.Cells.Item(125, 9).Value = 905
$ws.Cells.Item(125, 10).Value = 2
$ws.Cells.Item(125, 11).Value = 1
$ws.Cells.Item(125, 12).Value = 9
$ws.Cells.Item(125, 13).Value = 2
$ws.Cells.Item(125, 14).Value = 2
$ws.Cells.Item(125, 15).Value = 8
$ws.Cells.Item(125, 16).Value = 4
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = 9
$ws.Cells.Item(125, 21).Value = 1
$ws.Cells.Item(125, 22).Value = 3
$ws.Cells.Item(125, 23).Value = 20
$ws.Cells.Item(125, 24).Value = 13
$ws.Cells.Item(125, 25).Value = 1
$ws.Cells.Item(125, 26).Value = "T_1"
$ws.Cells.Item(125, 27).Value = "T_1"

# --- row 126 ---
$ws.Cells.Item(126, 1).Value = 36
$ws.Cells.Item(126, 2).Value = 2021
$ws.Cells.Item(126, 3).Value = 14
$ws.Cells.Item(126, 4).Value = 11
$ws.Cells.Item(126, 5).Value = 24
$ws.Cells.Item(126, 6).Value = "elabe"
$ws.Cells.Item(126, 7).Value = "online"
$ws.Cells.Item(126, 8).Value = "partially"
$ws.Cells.Item(126, 9).Value = 891
$ws.Cells.Item(126, 10).Value = 1
$ws.Cells.Item(126, 11).Value = 1
$ws.Cells.Item(126, 12).Value = 9
$ws.Cells.Item(126, 13).Value = 2
$ws.Cells.Item(126, 14).Value = 2
$ws.Cells.Item(126, 15).Value = 8
$ws.Cells.Item(126, 16).Value = 4
$ws.Cells.Item(126, 17).Value = 25
$ws.Cells.Item(126, 19).Value = 9
$ws.Cells.Item(126, 21).Value = 1
$ws.Cells.Item(126, 22).Value = 3
$ws.Cells.Item(126, 23).Value = 20
$ws.Cells.Item(126, 24).Value = 13
$ws.Cells.Item(126, 25).Value = 1
$ws.Cells.Item(126, 26).Value = "T_1"
$ws.Cells.Item(126, 27).Value = 1

# --- row 127 ---
$ws.Cells.Item(127, 1).Value = 37
$ws.Cells.Item(127, 2).Value = 2021
$ws.Cells.Item(127, 3).Value = 6
$ws.Cells.Item(127, 4).Value = 11
$ws.Cells.Item(127, 5).Value = 11
$ws.Cells.Item(127, 6).Value = "cluster17"
$ws.Cells.Item(127, 7).Value = "online"
$ws.Cells.Item(127, 8).Value = "partially"
$ws.Cells.Item(127, 9).Value = 6656
$ws.Cells.Item(127, 10).Value = 1
$ws.Cells.Item(127, 11).Value = 0.5
$ws.Cells.Item(127, 12).Value = 12
$ws.Cells.Item(127, 13).Value = 2
$ws.Cells.Item(127, 14).Value = 3
$ws.Cells.Item(127, 15).Value = 7
$ws.Cells.Item(127, 16).Value = 4
$ws.Cells.Item(127, 17).Value = 24
$ws.Cells.Item(127, 20).Value = 10
$ws.Cells.Item(127, 21).Value = 1
$ws.Cells.Item(127, 22).Value = 2
$ws.Cells.Item(127, 23).Value = 13
$ws.Cells.Item(127, 24).Value = 18.5
$ws.Cells.Item(127, 27).Value = 1.5
$ws.Cells.Item(127, 28).Value = 0.5

# --- row 128 ---
$ws.Cells.Item(128, 1).Value = 38
$ws.Cells.Item(128, 2).Value = 2021
$ws.Cells.Item(128, 3).Value = 11
$ws.Cells.Item(128, 4).Value = 11
$ws.Cells.Item(128, 5).Value = 11
$ws.Cells.Item(128, 6).Value = "cluster17"
$ws.Cells.Item(128, 7).Value = "online"
$ws.Cells.Item(128, 8).Value = "partially"
$ws.Cells.Item(128, 9).Value = 1727
$ws.Cells.Item(128, 10).Value = 1
$ws.Cells.Item(128, 11).Value = 0.5
$ws.Cells.Item(128, 12).Value = 11
$ws.Cells.Item(128, 13).Value = 1.5
$ws.Cells.Item(128, 14).Value = 3
$ws.Cells.Item(128, 15).Value = 6
$ws.Cells.Item(128, 16).Value = 4
$ws.Cells.Item(128, 17).Value = 23
$ws.Cells.Item(128, 20).Value = 11
$ws.Cells.Item(128, 21).Value = 1
$ws.Cells.Item(128, 22).Value = 3
$ws.Cells.Item(128, 23).Value = 14
$ws.Cells.Item(128, 24).Value = 19
$ws.Cells.Item(128, 27).Value = 2

# --- row 129 ---
$ws.Cells.Item(129, 1).Value = 39
$ws.Cells.Item(129, 2).Value = 2021
$ws.Cells.Item(129, 3).Value = 13
$ws.Cells.Item(129, 4).Value = 11
$ws.Cells.Item(129, 5).Value = 19
$ws.Cells.Item(129, 6).Value = "cluster17"
$ws.Cells.Item(129, 7).Value = "online"
$ws.Cells.Item(129, 8).Value = "partially"
$ws.Cells.Item(129, 9).Value = 1654
$ws.Cells.Item(129, 10).Value = 2
$ws.Cells.Item(129, 11).Value = 0.5
$ws.Cells.Item(129, 12).Value = 12
$ws.Cells.Item(129, 13).Value = 2.5
$ws.Cells.Item(129, 14).Value = 2
$ws.Cells.Item(129, 15).Value = 7
$ws.Cells.Item(129, 16).Value = 5
$ws.Cells.Item(129, 17).Value = 24
$ws.Cells.Item(129, 20).Value = 10
$ws.Cells.Item(129, 21).Value = 1
$ws.Cells.Item(129, 22).Value = 2
$ws.Cells.Item(129, 23).Value = 14
$ws.Cells.Item(129, 24).Value = 17
$ws.Cells.Item(129, 27).Value = 1

# --- row 130 ---
$ws.Cells.Item(130, 1).Value = 39
$ws.Cells.Item(130, 2).Value = 2021
$ws.Cells.Item(130, 3).Value = 13
$ws.Cells.Item(130, 4).Value = 11
$ws.Cells.Item(130, 5).Value = 19
$ws.Cells.Item(130, 6).Value = "cluster17"
$ws.Cells.Item(130, 7).Value = "online"
$ws.Cells.Item(130, 8).Value = "partially"
$ws.Cells.Item(130, 9).Value = 1630
$ws.Cells.Item(130, 10).Value = 2
$ws.Cells.Item(130, 11).Value = 0.5
$ws.Cells.Item(130, 12).Value = 12
$ws.Cells.Item(130, 13).Value = 2.5
$ws.Cells.Item(130, 14).Value = 2
$ws.Cells.Item(130, 15).Value = 8
$ws.Cells.Item(130, 16).Value = 4
$ws.Cells.Item(130, 17).Value = 24
$ws.Cells.Item(130, 18).Value = 9
$ws.Cells.Item(130, 21).Value = 1
$ws.Cells.Item(130, 22).Value = 3
$ws.Cells.Item(130, 23).Value = 15
$ws.Cells.Item(130, 24).Value = 16
$ws.Cells.Item(130, 27).Value = 1

# --- row 131 ---
$ws.Cells.Item(131, 1).Value = 39
$ws.Cells.Item(131, 2).Value = 2021
$ws.Cells.Item(131, 3).Value = 13
$ws.Cells.Item(131, 4).Value = 11
$ws.Cells.Item(131, 5).Value = 19
$ws.Cells.Item(131, 6).Value = "cluster17"
$ws.Cells.Item(131, 7).Value = "online"
$ws.Cells.Item(131, 8).Value = "partially"
$ws.Cells.Item(131, 9).Value = 1581
$ws.Cells.Item(131, 10).Value = 2
$ws.Cells.Item(131, 11).Value = 0.5
$ws.Cells.Item(131, 12).Value = 12
$ws.Cells.Item(131, 13).Value = 2.5
$ws.Cells.Item(131, 14).Value = 2
$ws.Cells.Item(131, 15).Value = 7
$ws.Cells.Item(131, 16).Value = 5
$ws.Cells.Item(131, 17).Value = 24
$ws.Cells.Item(131, 19).Value = 8
$ws.Cells.Item(131, 21).Value = 1
$ws.Cells.Item(131, 22).Value = 3
$ws.Cells.Item(131, 23).Value = 14
$ws.Cells.Item(131, 24).Value = 17
$ws.Cells.Item(131, 27).Value = 2

# --- Reproduce the final selection state (bottom-right corner of the new
#     used range) recorded in the workbook after the edit ----------------
$ws.Range("AB131").Select()
